$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be stored as TEXT (not auto-converted
# to a number) by round-tripping it through a blank staging cell (A1, which is
# inside the existing used range so it leaves no structural trace) formatted as
# Text, then pasting values+formats into the destination and clearing the staging
# cell completely so no residual style/content remains.
function Set-TextValue($cellRef, $text) {
    $ws.Range("A1").NumberFormat = "@"
    $ws.Range("A1").Value = $text
    $ws.Range("A1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4104)
    $ws.Range("A1").Clear()
}

$ws.Range('D2').Value = '67.512.63'
$ws.Range('E2').Value = '  -1.91%  '
$ws.Range('D3').Value = '3.511.88'
$ws.Range('E3').Value = '  -4.44%  '
$ws.Range('E4').Value = '  +0.16%  '
Set-TextValue 'D5' '197.88'
$ws.Range('E5').Value = '  -2.63%  '
Set-TextValue 'D6' '551.17'
$ws.Range('E6').Value = '  -5.31%  '
Set-TextValue 'D7' '0.637'
$ws.Range('E7').Value = '  +2.11%  '
$ws.Range('D8').Value = '3.500.30'
$ws.Range('E8').Value = '  -4.69%  '
$ws.Range('E9').Value = '  +0.06%  '
Set-TextValue 'D10' '0.653'
$ws.Range('E10').Value = '  -5.01%  '
Set-TextValue 'D11' '61.05'
$ws.Range('E11').Value = '  +5.36%  '
$ws.Range('E12').Value = '  -9.03%  '
Set-TextValue 'D13' '0.0000267'
$ws.Range('E13').Value = '  -10.62%  '
Set-TextValue 'D14' '9.76'
$ws.Range('E14').Value = '  -4.68%  '
$ws.Range('D15').Value = '4.076.21'
$ws.Range('E15').Value = '  -4.28%  '
$ws.Range('D16').Value = '3.509.39'
$ws.Range('E16').Value = '  -4.32%  '
$ws.Range('E17').Value = '  -2.07%  '
$ws.Range('D18').Value = '67.269.53'
$ws.Range('E18').Value = '  -2.18%  '
Set-TextValue 'D19' '18.29'
$ws.Range('E19').Value = '  -2.72%  '
Set-TextValue 'D20' '11.78'
$ws.Range('E20').Value = '  -6.93%  '
$ws.Range('E21').Value = '  -7.01%  '
Set-TextValue 'D22' '392.80'
$ws.Range('E22').Value = '  -3.60%  '
Set-TextValue 'D23' '3.96'
$ws.Range('E23').Value = '  -7.95%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D24' '85.32'
$ws.Range('E24').Value = '  -1.50%  '
$ws.Range('B25').Value = 'RenderToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D25' '11.76'
$ws.Range('E25').Value = '  -10.26%  '
Set-TextValue 'D26' '3.84'
$ws.Range('E26').Value = '  -1.62%  '
Set-TextValue 'D27' '12.22'
$ws.Range('E27').Value = '  -4.12%  '
Set-TextValue 'D28' '2.80'
$ws.Range('E28').Value = '  -6.37%  '
Set-TextValue 'D29' '8.83'
$ws.Range('E29').Value = '  -5.29%  '
Set-TextValue 'D30' '720.54'
$ws.Range('E30').Value = '  +3.36%  '
Set-TextValue 'D31' '31.19'
$ws.Range('E31').Value = '  -3.59%  '
Set-TextValue 'D32' '6.94'
$ws.Range('E32').Value = '  -16.68%  '
Set-TextValue 'D33' '11.66'
$ws.Range('E33').Value = '  -5.78%  '
Set-TextValue 'D34' '63.93'
$ws.Range('E34').Value = '  -1.79%  '
Set-TextValue 'D35' '0.110'
$ws.Range('E35').Value = '  -6.55%  '
Set-TextValue 'D36' '38.33'
$ws.Range('E36').Value = '  -11.68%  '
$ws.Range('E37').Value = '  +0.01%  '
Set-TextValue 'D38' '0.391'
$ws.Range('E38').Value = '  -9.50%  '
Set-TextValue 'D39' '3.00'
$ws.Range('E39').Value = '  -5.65%  '
$ws.Range('E40').Value = '  -8.65%  '
$ws.Range('E41').Value = '  +0.10%  '
$ws.Range('D42').Value = '3.059.18'
$ws.Range('E42').Value = '  -5.81%  '
$ws.Range('D43').Value = '0.0₃0676'
$ws.Range('E43').Value = '  -16.71%  '
$ws.Range('E44').Value = '  +5.19%  '
Set-TextValue 'D45' '2.49'
$ws.Range('E45').Value = '  -13.53%  '
Set-TextValue 'D46' '0.0406'
$ws.Range('E46').Value = '  -4.86%  '
Set-TextValue 'D47' '0.131'
$ws.Range('E47').Value = '  -1.65%  '
Set-TextValue 'D48' '2.56'
$ws.Range('E48').Value = '  -15.96%  '
Set-TextValue 'D49' '138.40'
$ws.Range('E49').Value = '  -3.20%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue 'D50' '2.91'
$ws.Range('E50').Value = '  -6.43%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue 'D51' '8.21'
$ws.Range('E51').Value = '  -9.14%  '
